# Update "Model Summary.xlsx" - Sheet1 with new model list using moving-average
# rice/wheat percentages (rice_moving_perc / wheat_moving_perc), replacing the
# old "as individual" / "as average" / "bpl change rate" models.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Clear the old region (A1:C10) completely, including the C column (which held
# the "bpl change rate" labels that are removed entirely in the new layout).
$ws.Range("A1:C10").Clear()

# Header row
$ws.Range("A1").Value = "Model"
$ws.Range("B1").Value = "R squared"

# Write the new model-name labels in the same order the workbook author
# originally entered them (moving-average variants first, then the plain
# "_perc" variants, then the last-added wheat moving model and the final
# bpl wheat "_perc" model) so the rebuilt shared-string table lines up.
$ws.Range("A3").Value = "Rice alt ~ Pop + rice_moving_perc + wheat_moving_perc"
$ws.Range("A6").Value = "Rice alt ~ bpl pop + rice_moving_perc + wheat_moving_perc"
$ws.Range("A13").Value = "Wheat alt ~ bpl pop + rice_moving_perc + wheat_moving_perc"
$ws.Range("A2").Value = "Rice Alt ~ Pop + rice_perc + wheat_perc "
$ws.Range("A5").Value = "Rice alt ~ bpl pop + rice_perc + wheat_perc "
$ws.Range("A9").Value = "Wheat alt ~ pop + rice_perc + wheat_perc "
$ws.Range("A10").Value = "Wheat alt ~ pop + rice_moving_perc + wheat_moving_perc"
$ws.Range("A12").Value = "Wheat alt ~ bpl pop + rice_perc + wheat_perc"

# R-squared values
$ws.Range("B2").Value = 0.77
$ws.Range("B3").Value = 0.781
$ws.Range("B5").Value = 0.6252
$ws.Range("B6").Value = 0.633
$ws.Range("B9").Value = 0.874
$ws.Range("B10").Value = 0.871
$ws.Range("B12").Value = 0.846
$ws.Range("B13").Value = 0.841

# Match the final selection state recorded in the saved workbook
$ws.Range("B13").Select()
